$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row from column A (header "Image" is in row 1,
# data rows follow in row 2..N).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "color/*") {
        $cell.Value = $val -replace "^color/", ""
    }
}
